# Adds more interview questions and an example for init/destroy methods
# (InitializingBean/DisposableBean, @Autowired, autowiring types,
#  @Component/@Controller/@Repository/@Service, PersistenceExceptionTranslationPostProcessor,
#  component scanning, @ComponentScan) to Sheet1, rows 22-28.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 22 ---
$ws.Range("A17").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("A22").Value2 = "What is InitializingBean and DisposableBean?"

# --- Row 23 ---
$ws.Range("A17").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("A23").Value2 = "What is @Autowired annotation?"

# --- Row 24 ---
$ws.Range("A17").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").Value2 = "What are different types of Spring Bean autowiring?"

$ws.Range("A17").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("B24").Value2 = "byName, byType and byConstructor"

# --- Row 25 (taller row, hyperlink in B) ---
$ws.Range("A17").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("A25").Value2 = "What" + [char]8217 + "s the difference between @Component, @Controller, @Repository & @Service annotations in Spring?"
$ws.Range("A25").EntireRow.RowHeight = 30

$ws.Hyperlinks.Add($ws.Range("B25"), "https://javarevisited.blogspot.com/2017/11/difference-between-component-service.html")
$ws.Range("B18").Copy()
$ws.Range("B25").PasteSpecial(-4122)

# --- Row 26 ---
$ws.Range("A17").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("A26").Value2 = "What is PersistenceExceptionTranslationPostProcessor ?"

# --- Row 27 (hyperlink in B, same url as row 25) ---
$ws.Range("A17").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value2 = "How does componenet scanning works in spring?"

$ws.Hyperlinks.Add($ws.Range("B27"), "https://javarevisited.blogspot.com/2017/11/difference-between-component-service.html")
$ws.Range("B18").Copy()
$ws.Range("B27").PasteSpecial(-4122)

# --- Row 28 ---
$ws.Range("A17").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value2 = "What is @ComponentScan?"

# Selection matches the end of the authored edit.
$ws.Range("B28").Select() | Out-Null
